# Add a short link to the updated resume in the document footer.
#
# The footer currently contains a single empty paragraph. We replace it
# with a paragraph that:
#   - uses the "Header & Footer" paragraph style (same style used
#     throughout the header),
#   - re-points the style's inherited right tab stop (9020) to a
#     center/right tab pair (4680 / 9360), matching the tab layout
#     already used in the document header,
#   - starts with a tab character to push the text to the center tab
#     stop,
#   - is followed by the "View the most recent version..." label and the
#     short link text, each carrying the rtl/lang run properties used
#     elsewhere in the document.

$d = $word.ActiveDocument
$footer = $d.Sections.Item(1).Footers.Item(1)

# Build the replacement paragraph as a WordprocessingML fragment so the
# paragraph formatting (including the cleared inherited tab stop) and run
# formatting can be expressed exactly, then drop it into the footer's
# range - this replaces the existing "<w:p><w:r/></w:p>" paragraph.
$footerParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Header &amp; Footer"/>
    <w:tabs>
      <w:tab w:val="center" w:pos="4680"/>
      <w:tab w:val="right" w:pos="9360"/>
      <w:tab w:val="clear" w:pos="9020"/>
    </w:tabs>
    <w:jc w:val="left"/>
  </w:pPr>
  <w:r>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rtl w:val="0"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">View the most recent version of this resume on GitHub: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rtl w:val="0"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>https://runty.link/resume</w:t>
  </w:r>
</w:p>
'@

$null = $footer.Range.InsertXML($footerParagraphXml)
